# Auto-generated edit script applying the cryptos.xlsx diff.
#
# All target cells are plain text (t="inlineStr") in the original workbook,
# including ones that look numeric (e.g. "0.999", "595.59"). A direct
# $range.Value = "0.999" would be auto-coerced to a real Excel number, so
# every write is prefixed with a literal leading apostrophe, which is the
# standard Excel "force text" input convention; Excel strips the apostrophe
# from the stored value but keeps the cell as Text. That quote-prefix does
# add a quotePrefix style flag, so we immediately reset Range.Style back to
# "Normal" afterwards to keep styling identical (unstyled, index 0) to the
# surrounding untouched cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.993.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +3.03%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.036.55'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +2.11%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  +0.05%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''595.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +1.27%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''152.53'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +7.10%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -0.05%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''3.031.34'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +2.02%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.521'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +0.63%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  +10.24%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.152'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +5.95%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.463'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +2.19%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.0000235'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +4.14%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''35.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +2.83%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = '''  +2.69%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''3.536.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +2.10%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''62.889.32'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +2.97%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  +0.55%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''3.034.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +2.20%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''455.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +1.47%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''14.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +2.34%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.694'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +1.65%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''7.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +2.11%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''82.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.90%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''2.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +5.14%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''10.85'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +9.55%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''12.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +0.83%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -0.02%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = '''NEARProtocol'
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = '''7.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +9.02%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = '''PancakeSwap'
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = '''2.71'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +2.49%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  +0.17%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  +5.59%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''27.66'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +1.89%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +4.41%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.0₃0854'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +7.82%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''1.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +2.41%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  +3.33%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''3.12'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +12.39%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''2.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +1.91%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''50.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +0.61%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''9.12'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -0.01%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  +4.43%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.294'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +11.69%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''41.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +10.62%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''392.01'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +0.47%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.0357'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +1.66%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''2.743.77'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +1.89%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''132.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +0.77%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = '''2.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +2.78%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  +0.81%  '
$ws.Range("E51").Style = "Normal"
